# ============================================================================
# MAJ automatique BRVM via GitHub Actions
# Refreshes the "Recommandations" and "Top_YTD" sheets with the latest values
# (rankings reshuffle as prices move, so row order/content both change).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---- Sheet "Recommandations": columns A..G, data rows 2..45 ----
$wsReco = $wb.Worksheets.Item("Recommandations")

$recoData = @(
    @("BRVM - SERVICES PUBLICS", 0, 8, 3258.43, 99.36, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI", 0, 4, 2915, 730, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI", 0, 4, 2685, 665, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 4, 2532.9, 636.69, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI", 0, 4, 2385, 615, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI", 0, 4, 2355, 560, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 4, 2320, 590, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 4, 1498.02, 376.95, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 4, 1482.51, 365.27, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 4, 1317.62, 329.4, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE", 0, 4, 1070.38, 269.33, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE", 0, 4, 883.45, 221.86, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL", 0, 4, 769.85, 192.4, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 4, 568.85, 143.89, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 4, 527.15, 131.24, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 4, 497.88, 123.95, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 4, 489.32, 121.82, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 4, 440.13, 108.29, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 430.8, 107.61, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 4, 373.6, 93.22, "🟡 Observer", "➖ Neutre"),
    @("SAFCA CI (SAFC)", 3, 0, 17.77, 3.27, "🟢 Achat", "✅ Renforcer"),
    @("UNIWAX CI (UNXC)", 3, 1, 11.88, 4.93, "🟢 Achat", "✅ Renforcer"),
    @("SUCRIVOIRE (SCRC)", 2, 0, 10.19, 2.95, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 2, 1, 8.27, 4.03, "🟡 Observer", "👀 À surveiller"),
    @("AIR LIQUIDE CI (SIVC)", 1, 0, 7.27, 7.27, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 2, 1, 6.53, -2.81, "🟡 Observer", "👀 À surveiller"),
    @("SICOR CI (SICC)", 1, 0, 3.79, 3.79, "🟡 Observer", "➖ Neutre"),
    @("SERVAIR ABIDJAN CI (ABJC)", 1, 0, 2.56, 2.56, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 1, 1, 0.47, 7.14, "🟡 Observer", "👀 À surveiller"),
    @("SICABLE CI (CABC)", 1, 1, 0.47, 3.16, "🟡 Observer", "👀 À surveiller"),
    @("ORAGROUP TOGO (ORGT)", 1, 1, 0.19, 4.39, "🟡 Observer", "👀 À surveiller"),
    @("VIVO ENERGY CI (SHEC)", 1, 1, 0.04, -2.13, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 0, 1, -0.82, -0.82, "🟡 Observer", "➖ Neutre"),
    @("NSIA BANQUE COTE D'IVOIRE (NSBC)", 0, 1, -0.89, -0.89, "🟡 Observer", "➖ Neutre"),
    @("SOLIBRA CI (SLBC)", 1, 1, -1.23, 3.83, "🟡 Observer", "👀 À surveiller"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -1.7, -1.7, "🟡 Observer", "➖ Neutre"),
    @("CIE CI (CIEC)", 0, 1, -1.86, -1.86, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 0, 1, -2.21, -2.21, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 0, 1, -2.38, -2.38, "🟡 Observer", "➖ Neutre"),
    @("BICI CI (BICC)", 0, 2, -2.4, -1.48, "🟡 Observer", "➖ Neutre"),
    @("ONATEL BF (ONTBF)", 0, 2, -2.91, -1.67, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 0, 1, -4.84, -4.84, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 0, 1, -5.88, -5.88, "🟡 Observer", "➖ Neutre"),
)

$r = 2
foreach ($row in $recoData) {
    $c = 1
    foreach ($val in $row) {
        $wsReco.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---- Sheet "Top_YTD": columns A..B, data rows 2..11 ----
$wsYtd = $wb.Worksheets.Item("Top_YTD")

$ytdData = @(
    @("BRVM - SERVICES PUBLICS", 6995594.73),
    @("UNIWAX CI", 471417.6),
    @("CFAO MOTORS CI", 353668.51),
    @("BRVM - AUTRES SECTEURS", 288928.84),
    @("SETAO CI", 234734.6),
    @("NEI-CEDA CI", 224663),
    @("AIR LIQUIDE CI", 213430.85),
    @("BRVM - TRANSPORT", 50593.09),
    @("BRVM - DISTRIBUTION", 48953.2),
    @("BRVM - AGRICULTURE", 33898.64),
)

$r = 2
foreach ($row in $ytdData) {
    $c = 1
    foreach ($val in $row) {
        $wsYtd.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
